# Bill of Materials update: the I2C level shifter part (row 11) was
# re-sourced from the PI6ULS5V9306 (Mouser 729-PI6ULS5V9306UEX) to the
# PCA9306DCUR (Mouser 595-PCA9306DCUR), at a lower unit price.
# Downstream formulas (row total I11, grand total J4) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E11").Value = "PCA9306DCUR"
$ws.Range("G11").Value = "595-PCA9306DCUR"
$ws.Range("H11").Value = 0.79
